$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 524.6667
$ws.Range("I2").Value = 310.21054
$ws.Range("J2").Value = 1339.6
$ws.Range("K2").Value = 310.21054
$ws.Range("L2").Value = 1339.6
$ws.Range("M2").Value = -197.21054
$ws.Range("N2").Value = -1565.6
$ws.Range("H17").Value = 3788.1875
$ws.Range("J17").Value = 3878.1292
$ws.Range("L17").Value = 11634.3876
$ws.Range("N17").Value = -11970.3876
$ws.Range("H40").Value = 4989.8
$ws.Range("J40").Value = 4989.8
$ws.Range("L40").Value = 4989.8
$ws.Range("N40").Value = -5339.8
$ws.Range("H57").Value = 92444.75
$ws.Range("J57").Value = 111259.664
$ws.Range("L57").Value = 333778.992
$ws.Range("N57").Value = -334776.992
$ws.Range("H69").Value = 7855
$ws.Range("I69").Value = 7318.75
$ws.Range("K69").Value = 21956.25
$ws.Range("M69").Value = -21082.25
$ws.Range("H72").Value = 7855
$ws.Range("I72").Value = 7318.75
$ws.Range("K72").Value = 65868.75
$ws.Range("M72").Value = -61500.75
$ws.Range("H112").Value = 3572.087
$ws.Range("J112").Value = 3583.121
$ws.Range("L112").Value = 10749.363
$ws.Range("N112").Value = -12965.363
$ws.Range("H138").Value = 119631.35
$ws.Range("I138").Value = 6054
$ws.Range("J138").Value = 137331.72
$ws.Range("K138").Value = 18162
$ws.Range("L138").Value = 411995.16
$ws.Range("M138").Value = -13022
$ws.Range("N138").Value = -422275.16

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22565.922
$ws.Range("I32").Value = 20214.486
$ws.Range("J32").Value = 49999.332
$ws.Range("K32").Value = 20214.486
$ws.Range("L32").Value = 49999.332
$ws.Range("M32").Value = -19927.486
$ws.Range("N32").Value = -50573.332
$ws.Range("H45").Value = 2698.8
$ws.Range("I45").Value = 2075.0588
$ws.Range("J45").Value = 4024.25
$ws.Range("K45").Value = 2075.0588
$ws.Range("L45").Value = 4024.25
$ws.Range("M45").Value = -1698.0588
$ws.Range("N45").Value = -4778.25
$ws.Range("H61").Value = 3591724
$ws.Range("I61").Value = 7144900
$ws.Range("K61").Value = 7144900
$ws.Range("M61").Value = -7144688
$ws.Range("H102").Value = 8003516
$ws.Range("I102").Value = 3709.9048
$ws.Range("K102").Value = 3709.9048
$ws.Range("M102").Value = -2087.9048
$ws.Range("H132").Value = 2082459.2
$ws.Range("I132").Value = 2655219.5
$ws.Range("K132").Value = 7965658.5
$ws.Range("M132").Value = -7963128.5
$ws.Range("H136").Value = 3591724
$ws.Range("I136").Value = 7144900
$ws.Range("K136").Value = 21434700
$ws.Range("M136").Value = -21432150

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3322.7917
$ws.Range("I99").Value = 2581.6316
$ws.Range("K99").Value = 2581.6316
$ws.Range("M99").Value = -1083.6316
$ws.Range("H134").Value = 3217.5173
$ws.Range("I134").Value = 2884.8262
$ws.Range("K134").Value = 8654.4786
$ws.Range("M134").Value = -6119.4786

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3998.6072
$ws.Range("I31").Value = 3737.611
$ws.Range("J31").Value = 4468.4
$ws.Range("K31").Value = 3737.611
$ws.Range("L31").Value = 4468.4
$ws.Range("M31").Value = -3442.611
$ws.Range("N31").Value = -5058.4
$ws.Range("H34").Value = 3998.6072
$ws.Range("I34").Value = 3737.611
$ws.Range("J34").Value = 4468.4
$ws.Range("K34").Value = 3737.611
$ws.Range("L34").Value = 4468.4
$ws.Range("M34").Value = -3535.611
$ws.Range("N34").Value = -4872.4
$ws.Range("H52").Value = 99165.836
$ws.Range("J52").Value = 108999
$ws.Range("L52").Value = 108999
$ws.Range("N52").Value = -109587
$ws.Range("H58").Value = 5086.4
$ws.Range("I58").Value = 6539.2
$ws.Range("K58").Value = 6539.2
$ws.Range("M58").Value = -6336.2
$ws.Range("H99").Value = 5999.7
$ws.Range("I99").Value = 5928.2856
$ws.Range("J99").Value = 6166.3335
$ws.Range("K99").Value = 5928.2856
$ws.Range("L99").Value = 6166.3335
$ws.Range("M99").Value = -4430.2856
$ws.Range("N99").Value = -9162.333500000001
$ws.Range("H126").Value = 5999.7
$ws.Range("I126").Value = 5928.2856
$ws.Range("J126").Value = 6166.3335
$ws.Range("K126").Value = 17784.8568
$ws.Range("L126").Value = 18499.0005
$ws.Range("M126").Value = -15314.8568
$ws.Range("N126").Value = -23439.0005
$ws.Range("H132").Value = 4997.7334
$ws.Range("I132").Value = 4101.7896
$ws.Range("K132").Value = 12305.3688
$ws.Range("M132").Value = -9775.3688
$ws.Range("H133").Value = 97250
$ws.Range("J133").Value = 97250
$ws.Range("L133").Value = 97250
$ws.Range("N133").Value = -102310
$ws.Range("H134").Value = 3962.0303
$ws.Range("I134").Value = 1726.2667
$ws.Range("K134").Value = 5178.800099999999
$ws.Range("M134").Value = -2643.800099999999
$ws.Range("H136").Value = 5086.4
$ws.Range("I136").Value = 6539.2
$ws.Range("K136").Value = 19617.6
$ws.Range("M136").Value = -17067.6
$ws.Range("H137").Value = 159500
$ws.Range("J137").Value = 159500
$ws.Range("L137").Value = 159500
$ws.Range("N137").Value = -169700
$ws.Range("H139").Value = 30909.092
$ws.Range("I139").Value = 30909.092
$ws.Range("K139").Value = 30909.092
$ws.Range("M139").Value = -25769.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2399.5
$ws.Range("I122").Value = 1466.6666
$ws.Range("K122").Value = 13199.9994
$ws.Range("M122").Value = -10749.9994
$ws.Range("H137").Value = 4891.5
$ws.Range("J137").Value = 6629.3335
$ws.Range("L137").Value = 19888.0005
$ws.Range("N137").Value = -30088.0005
$ws.Range("H139").Value = 4004145.2
$ws.Range("I139").Value = 2530.1333
$ws.Range("K139").Value = 7590.3999
$ws.Range("M139").Value = -2450.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 21013.25
$ws.Range("J24").Value = 26851.334
$ws.Range("L24").Value = 26851.334
$ws.Range("N24").Value = -27197.334
$ws.Range("H134").Value = 68581.25
$ws.Range("J134").Value = 68581.25
$ws.Range("L134").Value = 205743.75
$ws.Range("N134").Value = -210813.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14326.353
$ws.Range("I7").Value = 16122.637
$ws.Range("K7").Value = 16122.637
$ws.Range("M7").Value = -16010.637
$ws.Range("H46").Value = 5694.486
$ws.Range("I46").Value = 3031.4814
$ws.Range("J46").Value = 14682.125
$ws.Range("K46").Value = 3031.4814
$ws.Range("L46").Value = 14682.125
$ws.Range("M46").Value = -2843.4814
$ws.Range("N46").Value = -15058.125
$ws.Range("H126").Value = 14326.353
$ws.Range("I126").Value = 16122.637
$ws.Range("K126").Value = 48367.911
$ws.Range("M126").Value = -45897.911
$ws.Range("H136").Value = 6064.25
$ws.Range("I136").Value = 5208.6523
$ws.Range("K136").Value = 15625.9569
$ws.Range("M136").Value = -13075.9569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2197.75
$ws.Range("I96").Value = 1889.7142
$ws.Range("K96").Value = 1889.7142
$ws.Range("M96").Value = -516.7141999999999
$ws.Range("H132").Value = 1988.8182
$ws.Range("I132").Value = 1930.7778
$ws.Range("K132").Value = 5792.3334
$ws.Range("M132").Value = -3262.3334
$ws.Range("H136").Value = 16751.732
$ws.Range("I136").Value = 27285.75
$ws.Range("K136").Value = 81857.25
$ws.Range("M136").Value = -79307.25
